$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old dialogue/investigate rows (2-16); ClearContents keeps the existing
# cell styles (column B wrap-text style, J:L highlight style) so we only need to
# restate values, not formatting.
$ws.Range("A2:P16").ClearContents()

# Row 2: Yao suspects the bookshelf
$ws.Range("A2").Value = "Yao"
$ws.Range("B2").Value = "This bookshelf is filled entirely with martial arts manuals."
$ws.Range("C2").Value = "Yao-Regular"
$ws.Range("D2").Value = "DialogueVocal"
$ws.Range("E2").Value = "StudyInvestigate"
$ws.Range("F2").Value = "Suspicious"
$ws.Range("J2").Value = "appearAt"
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = "Dee-Thinking"

# Row 3: Dee on the Lord's skill
$ws.Range("A3").Value = "Dee"
$ws.Range("B3").Value = "I suppose the Lord’s formidable skills came from diligently studying these books."
$ws.Range("C3").Value = "Dee-Determined"
$ws.Range("D3").Value = "DialogueVocal"
$ws.Range("E3").Value = "StudyInvestigate"

# Row 4: Dee notices a particular book
$ws.Range("A4").Value = "Dee"
$ws.Range("B4").Value = "Hm? What’s this one?"
$ws.Range("C4").Value = "Dee-Thinking2"
$ws.Range("D4").Value = "DialogueVocal"
$ws.Range("E4").Value = "StudyInvestigate"

# Row 5: He identifies Twin Blade Throw manual
$ws.Range("A5").Value = "He"
$ws.Range("B5").Value = "That book is Qingliu Manor’s exclusive martial arts manual. It’s called Twin Blade Throw."
$ws.Range("C5").Value = "He-Regular"
$ws.Range("D5").Value = "DialogueVocal"
$ws.Range("E5").Value = "StudyInvestigate"

# Row 6: He - signature technique
$ws.Range("A6").Value = "He"
$ws.Range("B6").Value = "It was the Lord’s signature technique."
$ws.Range("C6").Value = "He-Sad"
$ws.Range("D6").Value = "DialogueVocal"
$ws.Range("E6").Value = "StudyInvestigate"

# Row 7: He - why only the Lord could use it
$ws.Range("A7").Value = "He"
$ws.Range("B7").Value = "That’s why no one else in the manor could use it."
$ws.Range("C7").Value = "He-Sad"
$ws.Range("D7").Value = "DialogueVocal"
$ws.Range("E7").Value = "StudyInvestigate"

# Row 8: He - intended for Young Lord Ming
$ws.Range("A8").Value = "He"
$ws.Range("B8").Value = "He had intended to pass it down to Young Lord Ming, but unfortunately......"
$ws.Range("C8").Value = "He-Regular"
$ws.Range("D8").Value = "DialogueVocal"
$ws.Range("E8").Value = "StudyInvestigate"

# Row 9: Dee condolences to Butler He
$ws.Range("A9").Value = "Dee"
$ws.Range("B9").Value = "My condolences, Butler He. Could you tell us more about this technique?"
$ws.Range("C9").Value = "Dee-Thinking2"
$ws.Range("D9").Value = "DialogueVocal"
$ws.Range("E9").Value = "StudyInvestigate"

# Row 10: He explains the technique power (part 1)
$ws.Range("A10").Value = "He"
$ws.Range("B10").Value = "Once mastered, it allows the user to throw twin blades with precision, capable of striking enemies over a hundred meters away. "
$ws.Range("C10").Value = "He-Sad"
$ws.Range("D10").Value = "DialogueVocal"
$ws.Range("E10").Value = "StudyInvestigate"

# Row 11: He explains the technique power (part 2)
$ws.Range("A11").Value = "He"
$ws.Range("B11").Value = "The power is immense."
$ws.Range("C11").Value = "He-Sad"
$ws.Range("D11").Value = "DialogueVocal"
$ws.Range("E11").Value = "StudyInvestigate"

# Row 12: Suspicious inner thought (colored)
$ws.Range("B12").Value = " <color=#00CC00>(With such a formidable martial skill, how could the Lord have been killed?)</color>"
$ws.Range("D12").Value = "DialogueVocal"
$ws.Range("E12").Value = "StudyInvestigate"

# Row 13: Investigate - Desk
$ws.Range("A13").Value = "Investigate"
$ws.Range("B13").Value = "Desk"
$ws.Range("C13").Value = "Desk"
$ws.Range("D13").Value = "DialogueVocal"
$ws.Range("E13").Value = "StudyInvestigate"

# Row 14: Investigate - Paper
$ws.Range("B14").Value = "Paper"
$ws.Range("C14").Value = "Paper"
$ws.Range("D14").Value = "DialogueVocal"
$ws.Range("E14").Value = "StudyInvestigate"

# Row 15: Investigate - Book
$ws.Range("B15").Value = "Book"
$ws.Range("C15").Value = "Book"
$ws.Range("D15").Value = "DialogueVocal"
$ws.Range("E15").Value = "StudyInvestigate"

# Row 16: End Investigation script trigger
$ws.Range("B16").Value = "End Investigation"
$ws.Range("C16").Value = "StoryScript13"
$ws.Range("D16").Value = "DialogueVocal"
$ws.Range("E16").Value = "StudyInvestigate"

# Row 17: disappear action
$ws.Range("J17").Value = "disappear"

# J16:L16 and K17:L17 need the same highlighted style as the rest of the J:L column
# (those cells did not exist in the old sheet, so ClearContents left no style behind).
$ws.Range("J15:L15").Copy()
$ws.Range("J16:L16").PasteSpecial(-4122)
$ws.Range("J15:L15").Copy()
$ws.Range("J17:L17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore explicit row heights where they differ from what carried over automatically
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 17
$ws.Rows.Item(9).RowHeight = 34
$ws.Rows.Item(10).RowHeight = 51
$ws.Rows.Item(12).RowHeight = 34
$ws.Rows.Item(15).RowHeight = 17
$ws.Rows.Item(16).RowHeight = 17

# Selection + dimension match the edited workbook
$ws.Range("B16").Select()
